$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.185348510742188
$ws.Range("B1").Value = 2.087687730789185
$ws.Range("C1").Value = 6.415035724639893
$ws.Range("D1").Value = 2.306306123733521
$ws.Range("E1").Value = 1.196385860443115
